$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Part 1: in-place value updates on rows 2-25
$ws.Range("F2").Value = $null
$ws.Range("F5").Value = 17.66
$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43
$ws.Range("E8").Value = $null
$ws.Range("F9").Value = $null
$ws.Range("F10").Value = $null
$ws.Range("E12").Value = -5.3
$ws.Range("E14").Value = $null
$ws.Range("E17").Value = -7.3
$ws.Range("E18").Value = -8.5
$ws.Range("E19").Value = $null
$ws.Range("E20").Value = $null
$ws.Range("E23").Value = -7
$ws.Range("F24").Value = 16.78

# Part 2: delete row 26 (RM 232) entirely
$ws.Rows(26).Delete()
# After that deletion, old row 28 (SC 92) is now row 27; delete it too
$ws.Rows(27).Delete()

# Part 3: apply remaining cell-level changes to the now-shifted rows
# row 27 (old SC 101, now at 27): C -> 10, E -> empty
$ws.Range("C27").Value = 10
$ws.Range("E27").Value = $null
# row 28 (old SC 105, now at 28): C -> empty, F -> empty
$ws.Range("C28").Value = $null
$ws.Range("F28").Value = $null
# row 29 (old SC 119, now at 29): C -> empty
$ws.Range("C29").Value = $null
# row 30 (old SC 120, now at 30): C -> 11.4, F -> 16.89
$ws.Range("C30").Value = 11.4
$ws.Range("F30").Value = 16.89
# row 32 (old SC 193, now at 32): C -> empty
$ws.Range("C32").Value = $null
